$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.479.70'
$ws.Range('D3').Value = '2.057.03'
$ws.Range('E3').Value = '  +3.89%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '252.77'
$ws.Range('E5').Value = '  +2.70%  '
$ws.Range('E6').Value = '  +2.88%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '68.16'
$ws.Range('E7').Value = '  +16.24%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  +6.47%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '59.61'
$ws.Range('E10').Value = '  +1.50%  '
$ws.Range('E11').Value = '  +4.56%  '
$ws.Range('E12').Value = '  +0.96%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.931'
$ws.Range('E13').Value = '  -0.63%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '14.88'
$ws.Range('E14').Value = '  +2.45%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '22.63'
$ws.Range('E15').Value = '  +25.04%  '
$ws.Range('D16').Value = '2.358.40'
$ws.Range('E16').Value = '  +3.86%  '
$ws.Range('E17').Value = '  +5.19%  '
$ws.Range('D18').Value = '2.062.71'
$ws.Range('E18').Value = '  +3.38%  '
$ws.Range('D19').Value = '37.378.30'
$ws.Range('E19').Value = '  +5.26%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '73.73'
$ws.Range('E20').Value = '  +3.26%  '
$ws.Range('D21').Value = '0.0₃0876'
$ws.Range('E21').Value = '  +3.58%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.50'
$ws.Range('E22').Value = '  +5.42%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '240.05'
$ws.Range('E23').Value = '  +3.18%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.71'
$ws.Range('E24').Value = '  +3.02%  '
$ws.Range('E26').Value = '  +5.95%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.95'
$ws.Range('E27').Value = '  +9.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '161.85'
$ws.Range('E28').Value = '  -1.81%  '
$ws.Range('E29').Value = '  +4.59%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.128'
$ws.Range('E30').Value = '  +31.96%  '
$ws.Range('E31').Value = '  +2.96%  '
$ws.Range('E32').Value = '  +7.94%  '
$ws.Range('E33').Value = '  +9.34%  '
$ws.Range('E34').Value = '  +5.71%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.67'
$ws.Range('E35').Value = '  +7.55%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.45'
$ws.Range('E36').Value = '  +1.34%  '
$ws.Range('B37').Value = 'THORChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.15'
$ws.Range('E37').Value = '  +14.43%  '
$ws.Range('B38').Value = 'BinanceUSD'
$ws.Range('C38').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  +0.11%  '
$ws.Range('E39').Value = '  +3.97%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.15'
$ws.Range('E40').Value = '  +38.73%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.103'
$ws.Range('E41').Value = '  +15.32%  '
$ws.Range('B42').Value = 'HuobiToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.04'
$ws.Range('E42').Value = '  +5.64%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.25'
$ws.Range('E43').Value = '  +2.04%  '
$ws.Range('E44').Value = '  +8.61%  '
$ws.Range('E45').Value = '  +5.72%  '
$ws.Range('E46').Value = '  +2.98%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '97.57'
$ws.Range('E47').Value = '  +4.14%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.94'
$ws.Range('E48').Value = '  +1.74%  '
$ws.Range('D49').Value = '1.411.75'
$ws.Range('E49').Value = '  +2.72%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.95'
$ws.Range('E50').Value = '  +2.01%  '
$ws.Range('E51').Value = '  +10.72%  '
